$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values (rows 2-25 correspond to data index 0-23).
# Columns F, H, I, J, L are unchanged (remain 0) and are not touched.
$colNumbers = @(2, 3, 4, 5, 7, 11, 13, 14, 15)

$col2Values = @(10.13519245628993, 9.879148322919058, 9.720958383948156, 9.656353297329114, 9.645620268173147, 9.720087528493391, 10.04717267927132, 10.67647706199697, 11.12600331166828, 11.32668854197044, 11.40205811038121, 11.38585492478071, 11.33290220425921, 11.30038345455596, 11.11280463162577, 10.9966984796084, 10.92956327140164, 10.90677412815501, 11.00909540432849, 11.34847327115344, 11.56659824085821, 11.45054156652136, 11.00349194544353, 10.50811458066004)
$col3Values = @(4.758500477031853, 4.58263332929997, 4.470119743864907, 4.423168367283535, 4.415306823394186, 4.469490946064184, 4.698823159373088, 5.111155798004337, 5.389669776109421, 5.510797533448574, 5.555846026265828, 5.546180769042863, 5.514520232490677, 5.495019939087668, 5.381640023301248, 5.310644129284494, 5.269285874031556, 5.255193430347829, 5.318256069413711, 5.523842091012483, 5.653418156319331, 5.584704503149453, 5.314816397344859, 5.00379509939122)
$col4Values = @(5.995016904908919, 5.878042332853814, 5.806809299308139, 5.777974211134683, 5.77319913129264, 5.806419578734295, 5.954588296736801, 6.247895321285775, 6.462549349015497, 6.559515806087353, 6.596097218077378, 6.588225490543917, 6.562528360474934, 6.546769051849777, 6.4561952148234, 6.400427765754794, 6.368289980585502, 6.357399307600092, 6.406370999427573, 6.570080273562182, 6.676254324440306, 6.619674797174158, 6.403684298222176, 6.168518619556606)
$col5Values = @(12.1428478763065, 11.92322238730697, 11.79045061519949, 11.73694725824284, 11.72810174648161, 11.7897265074149, 12.0667372427582, 12.62295678013585, 13.03491072452392, 13.2220891048106, 13.29286148379976, 13.27762536681526, 13.22791418018298, 13.1974483727211, 13.02266713550346, 12.91533257054175, 12.85358008063536, 12.83267128999722, 12.92676078656067, 13.24251908263349, 13.44822289352125, 13.33851980436762, 12.92159422263973, 12.47157603811033)
$col7Values = @(3.652230012369585, 3.654334734380805, 3.655695490515261, 3.656267277558416, 3.656363266945607, 3.655703131843975, 3.652941547829244, 3.648066641703836, 3.64481100634433, 3.643399951088615, 3.642875622063705, 3.642988101456778, 3.643356613964138, 3.643583640191699, 3.644904625261148, 3.645732885683918, 3.646215866224357, 3.64638052805015, 3.645644034588157, 3.643248101731476, 3.641740525713587, 3.642539829510363, 3.64568418297441, 3.64932793434914)
$col11Values = @(9.454700929734088, 9.286195383706158, 9.183274447216803, 9.141527830186451, 9.134609292334938, 9.182710573025144, 9.396522416906162, 9.817590355241872, 10.12475119710934, 10.26334853113489, 10.31561693391719, 10.30437038906672, 10.26765329346657, 10.2451334419069, 10.11566620269833, 10.03591411853283, 9.989939305225629, 9.974356889579889, 10.04441496801871, 10.27844424722947, 10.43011783086962, 10.34930045069002, 10.04057211701333, 9.703834891320335)
$col13Values = @(14.08188274895984, 13.9216510842417, 13.82599991482625, 13.7877525132647, 13.78144699016858, 13.82548107889402, 14.02609493841267, 14.43915999714971, 14.75186252390875, 14.8955253989997, 14.9500801466743, 14.93832480981983, 14.90001087682455, 14.87656088944241, 14.74249792538242, 14.66058048536824, 14.61359972902387, 14.59771775325792, 14.66928702943311, 14.91126086172781, 15.07027206542959, 14.98534193940699, 14.66535043981612, 14.32559289121147)
$col14Values = @(19.5942677194054, 19.65659305317522, 19.69668222208373, 19.71347800461593, 19.71629469589326, 19.69690687537639, 19.615380222497, 19.46990216304938, 19.37171733991818, 19.32892262710209, 19.3129850945335, 19.31640562857028, 19.32760607342653, 19.33450152744575, 19.3745516182113, 19.39959929068424, 19.41418212757946, 19.41914989838887, 19.39691470713849, 19.32430896865447, 19.27841796611623, 19.30276835376407, 19.39812783846623, 19.50772481070196)
$col15Values = @(24.96810809532453, 25.00055198984928, 25.02588845507682, 25.03757258777134, 25.03959474626331, 25.02604053125065, 24.97816930928805, 24.9273625578363, 24.91639876611451, 24.91714911306163, 24.91825850367742, 24.91798287648297, 24.91722384504167, 24.91686638391389, 24.91646521460483, 24.91768898701025, 24.91893296428978, 24.91944690213655, 24.91750281359626, 24.91742439530128, 24.92218303090479, 24.91920323880385, 24.91758529920338, 24.93648325742182)

$allValues = @(
    $col2Values, $col3Values, $col4Values, $col5Values, $col7Values, $col11Values, $col13Values, $col14Values, $col15Values
)

for ($c = 0; $c -lt $colNumbers.Length; $c++) {
    $colNum = $colNumbers[$c]
    $values = $allValues[$c]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $rowNum = 2 + $i
        $ws.Cells.Item($rowNum, $colNum).Value = $values[$i]
    }
}
